# Auto-generated Excel COM-interop script applying the ADD (Auxilio
# Desenvolvimento de Dados) stock-consistency refresh described in the
# commit "atualizei dados da add".
#
# It updates, per affected product/store row:
#   - C (estoque)          -> new stock count
#   - D (data_estoque)     -> new stock-check timestamp
#   - E (estoque_depois)   -> kept identical to C
#   - F (data_movimento)   -> new movement timestamp
# Three rows only touch F (the movement timestamp was recomputed while the
# stock counts stayed the same). One existing row (2570) gains a
# previously-missing F value, and three brand-new rows (2573-2575) are
# appended with freshly refreshed data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 45825.28852326881
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 45824.75037037037
$ws.Range("C5").Value = 14
$ws.Range("D5").Value = 45825.28852326892
$ws.Range("E5").Value = 14
$ws.Range("F5").Value = 45824.75209490741
$ws.Range("C7").Value = -2
$ws.Range("D7").Value = 45825.28852326469
$ws.Range("E7").Value = -2
$ws.Range("F7").Value = 45824.68784722222
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 45825.28849841389
$ws.Range("E24").Value = 20
$ws.Range("F24").Value = 45824.44024305556
$ws.Range("C33").Value = 3025
$ws.Range("D33").Value = 45825.2885232613
$ws.Range("E33").Value = 3025
$ws.Range("F33").Value = 45824.64774305555
$ws.Range("C49").Value = 66
$ws.Range("D49").Value = 45825.2884984123
$ws.Range("E49").Value = 66
$ws.Range("F49").Value = 45824.43540509259
$ws.Range("C56").Value = 351
$ws.Range("D56").Value = 45825.28849841777
$ws.Range("E56").Value = 351
$ws.Range("F56").Value = 45824.64238425926
$ws.Range("C58").Value = 161
$ws.Range("D58").Value = 45825.28852326544
$ws.Range("E58").Value = 161
$ws.Range("F58").Value = 45824.69262731481
$ws.Range("C59").Value = 140
$ws.Range("D59").Value = 45825.28852326087
$ws.Range("E59").Value = 140
$ws.Range("F59").Value = 45824.64711805555
$ws.Range("C63").Value = 196
$ws.Range("D63").Value = 45825.28849841583
$ws.Range("E63").Value = 196
$ws.Range("F63").Value = 45824.60457175926
$ws.Range("C70").Value = 75
$ws.Range("D70").Value = 45825.28852326251
$ws.Range("E70").Value = 75
$ws.Range("F70").Value = 45824.64796296296
$ws.Range("C81").Value = 9
$ws.Range("D81").Value = 45825.28849841786
$ws.Range("E81").Value = 9
$ws.Range("F81").Value = 45824.64238425926
$ws.Range("C82").Value = 17
$ws.Range("D82").Value = 45825.28852326557
$ws.Range("E82").Value = 17
$ws.Range("F82").Value = 45824.69262731481
$ws.Range("C83").Value = 12
$ws.Range("D83").Value = 45825.28849841796
$ws.Range("E83").Value = 12
$ws.Range("F83").Value = 45824.64238425926
$ws.Range("C85").Value = 73
$ws.Range("D85").Value = 45825.28849841443
$ws.Range("E85").Value = 73
$ws.Range("F85").Value = 45824.45010416667
$ws.Range("F86").Value = 45824.60457175926
$ws.Range("C87").Value = -15
$ws.Range("D87").Value = 45825.28852326569
$ws.Range("E87").Value = -15
$ws.Range("F87").Value = 45824.69262731481
$ws.Range("C91").Value = 96
$ws.Range("D91").Value = 45825.28849841641
$ws.Range("E91").Value = 96
$ws.Range("F91").Value = 45824.62041666666
$ws.Range("C94").Value = 87
$ws.Range("D94").Value = 45825.28849841835
$ws.Range("E94").Value = 87
$ws.Range("F94").Value = 45824.64261574074
$ws.Range("C96").Value = 5
$ws.Range("D96").Value = 45825.28849841553
$ws.Range("E96").Value = 5
$ws.Range("F96").Value = 45824.50672453704
$ws.Range("C101").Value = -9
$ws.Range("D101").Value = 45825.28852326581
$ws.Range("E101").Value = -9
$ws.Range("F101").Value = 45824.69262731481
$ws.Range("C106").Value = 28
$ws.Range("D106").Value = 45825.28852326265
$ws.Range("E106").Value = 28
$ws.Range("F106").Value = 45824.64796296296
$ws.Range("C115").Value = 204
$ws.Range("D115").Value = 45825.28849841651
$ws.Range("E115").Value = 204
$ws.Range("F115").Value = 45824.62041666666
$ws.Range("C117").Value = 1057
$ws.Range("D117").Value = 45825.2884984171
$ws.Range("E117").Value = 1057
$ws.Range("F117").Value = 45824.63096064814
$ws.Range("C120").Value = 7
$ws.Range("D120").Value = 45825.28852326947
$ws.Range("E120").Value = 7
$ws.Range("F120").Value = 45824.76686342592
$ws.Range("C123").Value = 379
$ws.Range("D123").Value = 45825.28852326279
$ws.Range("E123").Value = 379
$ws.Range("F123").Value = 45824.64796296296
$ws.Range("C139").Value = 35
$ws.Range("D139").Value = 45825.28849841905
$ws.Range("E139").Value = 35
$ws.Range("F139").Value = 45824.64331018519
$ws.Range("C150").Value = 2
$ws.Range("D150").Value = 45825.28849841242
$ws.Range("E150").Value = 2
$ws.Range("F150").Value = 45824.43540509259
$ws.Range("F162").Value = 45824.60457175926
$ws.Range("F175").Value = 45824.60457175926
$ws.Range("C185").Value = 133
$ws.Range("D185").Value = 45825.28849841207
$ws.Range("E185").Value = 133
$ws.Range("F185").Value = 45824.38501157407
$ws.Range("C195").Value = 1
$ws.Range("D195").Value = 45825.28849841134
$ws.Range("E195").Value = 1
$ws.Range("F195").Value = 45824.3782175926
$ws.Range("C217").Value = 99
$ws.Range("D217").Value = 45825.28849841253
$ws.Range("E217").Value = 99
$ws.Range("F217").Value = 45824.43540509259
$ws.Range("C222").Value = 6
$ws.Range("D222").Value = 45825.28849841264
$ws.Range("E222").Value = 6
$ws.Range("F222").Value = 45824.43540509259
$ws.Range("C236").Value = -3
$ws.Range("D236").Value = 45825.28849841845
$ws.Range("E236").Value = -3
$ws.Range("F236").Value = 45824.64261574074
$ws.Range("C247").Value = 926
$ws.Range("D247").Value = 45825.28852326146
$ws.Range("E247").Value = 926
$ws.Range("F247").Value = 45824.64774305555
$ws.Range("C255").Value = 0
$ws.Range("D255").Value = 45825.2884984141
$ws.Range("E255").Value = 0
$ws.Range("F255").Value = 45824.44436342592
$ws.Range("C270").Value = -4
$ws.Range("D270").Value = 45825.28849841147
$ws.Range("E270").Value = -4
$ws.Range("F270").Value = 45824.3782175926
$ws.Range("C272").Value = 301
$ws.Range("D272").Value = 45825.28852326293
$ws.Range("E272").Value = 301
$ws.Range("F272").Value = 45824.64796296296
$ws.Range("C274").Value = 383
$ws.Range("D274").Value = 45825.28849841854
$ws.Range("E274").Value = 383
$ws.Range("F274").Value = 45824.64261574074
$ws.Range("C283").Value = 200
$ws.Range("D283").Value = 45825.28852326739
$ws.Range("E283").Value = 200
$ws.Range("F283").Value = 45824.74089120371
$ws.Range("C287").Value = 502
$ws.Range("D287").Value = 45825.28849841275
$ws.Range("E287").Value = 502
$ws.Range("F287").Value = 45824.43540509259
$ws.Range("C291").Value = 572
$ws.Range("D291").Value = 45825.28852326162
$ws.Range("E291").Value = 572
$ws.Range("F291").Value = 45824.64774305555
$ws.Range("C292").Value = 4
$ws.Range("D292").Value = 45825.28849841865
$ws.Range("E292").Value = 4
$ws.Range("F292").Value = 45824.64261574074
$ws.Range("C295").Value = 138
$ws.Range("D295").Value = 45825.28852326792
$ws.Range("E295").Value = 138
$ws.Range("F295").Value = 45824.74362268519
$ws.Range("C309").Value = 990
$ws.Range("D309").Value = 45825.28849841504
$ws.Range("E309").Value = 990
$ws.Range("F309").Value = 45824.4725
$ws.Range("C314").Value = 6
$ws.Range("D314").Value = 45825.28852326102
$ws.Range("E314").Value = 6
$ws.Range("F314").Value = 45824.64711805555
$ws.Range("C315").Value = 68
$ws.Range("D315").Value = 45825.28849841914
$ws.Range("E315").Value = 68
$ws.Range("F315").Value = 45824.64331018519
$ws.Range("C318").Value = -134
$ws.Range("D318").Value = 45825.28852326493
$ws.Range("E318").Value = -134
$ws.Range("F318").Value = 45824.69230324074
$ws.Range("C319").Value = -1
$ws.Range("D319").Value = 45825.28849841288
$ws.Range("E319").Value = -1
$ws.Range("F319").Value = 45824.43540509259
$ws.Range("C326").Value = 18
$ws.Range("D326").Value = 45825.28849841631
$ws.Range("E326").Value = 18
$ws.Range("F326").Value = 45824.61862268519
$ws.Range("C328").Value = 109
$ws.Range("D328").Value = 45825.28849841077
$ws.Range("E328").Value = 109
$ws.Range("F328").Value = 45824.37760416666
$ws.Range("C342").Value = 158
$ws.Range("D342").Value = 45825.28849841953
$ws.Range("E342").Value = 158
$ws.Range("F342").Value = 45824.64376157407
$ws.Range("C343").Value = 14
$ws.Range("D343").Value = 45825.28849841874
$ws.Range("E343").Value = 14
$ws.Range("F343").Value = 45824.64261574074
$ws.Range("C346").Value = 133
$ws.Range("D346").Value = 45825.28849841515
$ws.Range("E346").Value = 133
$ws.Range("F346").Value = 45824.4725
$ws.Range("C350").Value = -270
$ws.Range("D350").Value = 45825.28852326308
$ws.Range("E350").Value = -270
$ws.Range("F350").Value = 45824.64796296296
$ws.Range("C351").Value = 1395
$ws.Range("D351").Value = 45825.28849841563
$ws.Range("E351").Value = 1395
$ws.Range("F351").Value = 45824.50672453704
$ws.Range("C361").Value = 29
$ws.Range("D361").Value = 45825.28852326638
$ws.Range("E361").Value = 29
$ws.Range("F361").Value = 45824.74033564814
$ws.Range("C363").Value = 579
$ws.Range("D363").Value = 45825.28852326674
$ws.Range("E363").Value = 579
$ws.Range("F363").Value = 45824.74057870371
$ws.Range("C379").Value = 22
$ws.Range("D379").Value = 45825.28849841159
$ws.Range("E379").Value = 22
$ws.Range("F379").Value = 45824.3782175926
$ws.Range("C386").Value = -1
$ws.Range("D386").Value = 45825.28849841095
$ws.Range("E386").Value = -1
$ws.Range("F386").Value = 45824.37760416666
$ws.Range("C388").Value = -7
$ws.Range("D388").Value = 45825.28852326594
$ws.Range("E388").Value = -7
$ws.Range("F388").Value = 45824.69262731481
$ws.Range("C389").Value = -5
$ws.Range("D389").Value = 45825.28852326606
$ws.Range("E389").Value = -5
$ws.Range("F389").Value = 45824.69262731481
$ws.Range("C390").Value = 178
$ws.Range("D390").Value = 45825.28852326684
$ws.Range("E390").Value = 178
$ws.Range("F390").Value = 45824.74057870371
$ws.Range("C394").Value = 187
$ws.Range("D394").Value = 45825.28849841923
$ws.Range("E394").Value = 187
$ws.Range("F394").Value = 45824.64331018519
$ws.Range("C410").Value = 1835
$ws.Range("D410").Value = 45825.28852326178
$ws.Range("E410").Value = 1835
$ws.Range("F410").Value = 45824.64774305555
$ws.Range("C414").Value = 66
$ws.Range("D414").Value = 45825.288498413
$ws.Range("E414").Value = 66
$ws.Range("F414").Value = 45824.43540509259
$ws.Range("C418").Value = 77
$ws.Range("D418").Value = 45825.28852326803
$ws.Range("E418").Value = 77
$ws.Range("F418").Value = 45824.74362268519
$ws.Range("C422").Value = 55
$ws.Range("D422").Value = 45825.2884984131
$ws.Range("E422").Value = 55
$ws.Range("F422").Value = 45824.43540509259
$ws.Range("C435").Value = 65
$ws.Range("D435").Value = 45825.28849841982
$ws.Range("E435").Value = 65
$ws.Range("F435").Value = 45824.64493055556
$ws.Range("C443").Value = 23
$ws.Range("D443").Value = 45825.28852326116
$ws.Range("E443").Value = 23
$ws.Range("F443").Value = 45824.64711805555
$ws.Range("C469").Value = -114
$ws.Range("D469").Value = 45825.28852326193
$ws.Range("E469").Value = -114
$ws.Range("F469").Value = 45824.64774305555
$ws.Range("C472").Value = 43
$ws.Range("D472").Value = 45825.28849841379
$ws.Range("E472").Value = 43
$ws.Range("F472").Value = 45824.43664351852
$ws.Range("C475").Value = 0
$ws.Range("D475").Value = 45825.28849841719
$ws.Range("E475").Value = 0
$ws.Range("F475").Value = 45824.63096064814
$ws.Range("C480").Value = 226
$ws.Range("D480").Value = 45825.28849841963
$ws.Range("E480").Value = 226
$ws.Range("F480").Value = 45824.64376157407
$ws.Range("C486").Value = 16
$ws.Range("D486").Value = 45825.28852326445
$ws.Range("E486").Value = 16
$ws.Range("F486").Value = 45824.67936342592
$ws.Range("C505").Value = 7
$ws.Range("D505").Value = 45825.28852326814
$ws.Range("E505").Value = 7
$ws.Range("F505").Value = 45824.74362268519
$ws.Range("C507").Value = 12
$ws.Range("D507").Value = 45825.28849841455
$ws.Range("E507").Value = 12
$ws.Range("F507").Value = 45824.45626157407
$ws.Range("C510").Value = 265
$ws.Range("D510").Value = 45825.28852326749
$ws.Range("E510").Value = 265
$ws.Range("F510").Value = 45824.74089120371
$ws.Range("C535").Value = 98
$ws.Range("D535").Value = 45825.28849841469
$ws.Range("E535").Value = 98
$ws.Range("F535").Value = 45824.45626157407
$ws.Range("C570").Value = 2610
$ws.Range("D570").Value = 45825.28849841806
$ws.Range("E570").Value = 2610
$ws.Range("F570").Value = 45824.64238425926
$ws.Range("C581").Value = 3
$ws.Range("D581").Value = 45825.28849841622
$ws.Range("E581").Value = 3
$ws.Range("F581").Value = 45824.60457175926
$ws.Range("C606").Value = 114
$ws.Range("D606").Value = 45825.28852326824
$ws.Range("E606").Value = 114
$ws.Range("F606").Value = 45824.74362268519
$ws.Range("C631").Value = 60
$ws.Range("D631").Value = 45825.2884984168
$ws.Range("E631").Value = 60
$ws.Range("F631").Value = 45824.62626157407
$ws.Range("C634").Value = 19
$ws.Range("D634").Value = 45825.28852326208
$ws.Range("E634").Value = 19
$ws.Range("F634").Value = 45824.64774305555
$ws.Range("C652").Value = 35
$ws.Range("D652").Value = 45825.28849841534
$ws.Range("E652").Value = 35
$ws.Range("F652").Value = 45824.47506944444
$ws.Range("C656").Value = 8
$ws.Range("D656").Value = 45825.28849841322
$ws.Range("E656").Value = 8
$ws.Range("F656").Value = 45824.43540509259
$ws.Range("C657").Value = 2031
$ws.Range("D657").Value = 45825.28852326759
$ws.Range("E657").Value = 2031
$ws.Range("F657").Value = 45824.74089120371
$ws.Range("C660").Value = 400
$ws.Range("D660").Value = 45825.28849841816
$ws.Range("E660").Value = 400
$ws.Range("F660").Value = 45824.64238425926
$ws.Range("C681").Value = -13
$ws.Range("D681").Value = 45825.28852326323
$ws.Range("E681").Value = -13
$ws.Range("F681").Value = 45824.64796296296
$ws.Range("C691").Value = 161
$ws.Range("D691").Value = 45825.2884984122
$ws.Range("E691").Value = 161
$ws.Range("F691").Value = 45824.38501157407
$ws.Range("C692").Value = -2
$ws.Range("D692").Value = 45825.28852326506
$ws.Range("E692").Value = -2
$ws.Range("F692").Value = 45824.69230324074
$ws.Range("C720").Value = 570
$ws.Range("D720").Value = 45825.28849841883
$ws.Range("E720").Value = 570
$ws.Range("F720").Value = 45824.6428125
$ws.Range("C726").Value = 0
$ws.Range("D726").Value = 45825.28852326433
$ws.Range("E726").Value = 0
$ws.Range("F726").Value = 45824.67894675926
$ws.Range("C729").Value = 18
$ws.Range("D729").Value = 45825.28852326518
$ws.Range("E729").Value = 18
$ws.Range("F729").Value = 45824.69230324074
$ws.Range("C731").Value = 78
$ws.Range("D731").Value = 45825.28852326224
$ws.Range("E731").Value = 78
$ws.Range("F731").Value = 45824.64774305555
$ws.Range("C732").Value = 269
$ws.Range("D732").Value = 45825.28852326337
$ws.Range("E732").Value = 269
$ws.Range("F732").Value = 45824.64796296296
$ws.Range("C735").Value = 46
$ws.Range("D735").Value = 45825.28849841993
$ws.Range("E735").Value = 46
$ws.Range("F735").Value = 45824.64493055556
$ws.Range("C741").Value = -27
$ws.Range("D741").Value = 45825.28852326352
$ws.Range("E741").Value = -27
$ws.Range("F741").Value = 45824.64796296296
$ws.Range("C770").Value = -2
$ws.Range("D770").Value = 45825.288498414
$ws.Range("E770").Value = -2
$ws.Range("F770").Value = 45824.44172453704
$ws.Range("C772").Value = 13
$ws.Range("D772").Value = 45825.28852326696
$ws.Range("E772").Value = 13
$ws.Range("F772").Value = 45824.74057870371
$ws.Range("C798").Value = 2
$ws.Range("D798").Value = 45825.28849842004
$ws.Range("E798").Value = 2
$ws.Range("F798").Value = 45824.64493055556
$ws.Range("C810").Value = -1
$ws.Range("D810").Value = 45825.28852326925
$ws.Range("E810").Value = -1
$ws.Range("F810").Value = 45824.76438657408
$ws.Range("C812").Value = 19
$ws.Range("D812").Value = 45825.28852326237
$ws.Range("E812").Value = 19
$ws.Range("F812").Value = 45824.64774305555
$ws.Range("C826").Value = 29
$ws.Range("D826").Value = 45825.28852326835
$ws.Range("E826").Value = 29
$ws.Range("F826").Value = 45824.74362268519
$ws.Range("C853").Value = 74
$ws.Range("D853").Value = 45825.28849841728
$ws.Range("E853").Value = 74
$ws.Range("F853").Value = 45824.63096064814
$ws.Range("C854").Value = 11
$ws.Range("D854").Value = 45825.28852326937
$ws.Range("E854").Value = 11
$ws.Range("F854").Value = 45824.76539351852
$ws.Range("C872").Value = 525
$ws.Range("D872").Value = 45825.28852326706
$ws.Range("E872").Value = 525
$ws.Range("F872").Value = 45824.74057870371
$ws.Range("C876").Value = 750
$ws.Range("D876").Value = 45825.28849842013
$ws.Range("E876").Value = 750
$ws.Range("F876").Value = 45824.64493055556
$ws.Range("C883").Value = 465
$ws.Range("D883").Value = 45825.28852326845
$ws.Range("E883").Value = 465
$ws.Range("F883").Value = 45824.74362268519
$ws.Range("C888").Value = 424
$ws.Range("D888").Value = 45825.28849841659
$ws.Range("E888").Value = 424
$ws.Range("F888").Value = 45824.62041666666
$ws.Range("C889").Value = 17
$ws.Range("D889").Value = 45825.28849841544
$ws.Range("E889").Value = 17
$ws.Range("F889").Value = 45824.50439814815
$ws.Range("C924").Value = 11
$ws.Range("D924").Value = 45825.28852326902
$ws.Range("E924").Value = 11
$ws.Range("F924").Value = 45824.7532175926
$ws.Range("C962").Value = 15
$ws.Range("D962").Value = 45825.28849841737
$ws.Range("E962").Value = 15
$ws.Range("F962").Value = 45824.63096064814
$ws.Range("C963").Value = 816
$ws.Range("D963").Value = 45825.28849842023
$ws.Range("E963").Value = 816
$ws.Range("F963").Value = 45824.64493055556
$ws.Range("C1017").Value = 369
$ws.Range("D1017").Value = 45825.28852326716
$ws.Range("E1017").Value = 369
$ws.Range("F1017").Value = 45824.74057870371
$ws.Range("C1024").Value = -11
$ws.Range("D1024").Value = 45825.28852326619
$ws.Range("E1024").Value = -11
$ws.Range("F1024").Value = 45824.69262731481
$ws.Range("C1025").Value = 11
$ws.Range("D1025").Value = 45825.28849841332
$ws.Range("E1025").Value = 11
$ws.Range("F1025").Value = 45824.43540509259
$ws.Range("C1039").Value = 89
$ws.Range("D1039").Value = 45825.28849841346
$ws.Range("E1039").Value = 89
$ws.Range("F1039").Value = 45824.43540509259
$ws.Range("C1062").Value = 580
$ws.Range("D1062").Value = 45825.28849841934
$ws.Range("E1062").Value = 580
$ws.Range("F1062").Value = 45824.64331018519
$ws.Range("C1110").Value = 138
$ws.Range("D1110").Value = 45825.28852326408
$ws.Range("E1110").Value = 138
$ws.Range("F1110").Value = 45824.65804398148
$ws.Range("C1133").Value = 77
$ws.Range("D1133").Value = 45825.2884984169
$ws.Range("E1133").Value = 77
$ws.Range("F1133").Value = 45824.62626157407
$ws.Range("C1147").Value = 54
$ws.Range("D1147").Value = 45825.28852326653
$ws.Range("E1147").Value = 54
$ws.Range("F1147").Value = 45824.74033564814
$ws.Range("C1149").Value = 15
$ws.Range("D1149").Value = 45825.2884984148
$ws.Range("E1149").Value = 15
$ws.Range("F1149").Value = 45824.45626157407
$ws.Range("C1150").Value = 32
$ws.Range("D1150").Value = 45825.28849841747
$ws.Range("E1150").Value = 32
$ws.Range("F1150").Value = 45824.63096064814
$ws.Range("C1156").Value = 23
$ws.Range("D1156").Value = 45825.2885232642
$ws.Range("E1156").Value = 23
$ws.Range("F1156").Value = 45824.67607638889
$ws.Range("C1195").Value = 12
$ws.Range("D1195").Value = 45825.28852326531
$ws.Range("E1195").Value = 12
$ws.Range("F1195").Value = 45824.69230324074
$ws.Range("C1196").Value = 6
$ws.Range("D1196").Value = 45825.2884984167
$ws.Range("E1196").Value = 6
$ws.Range("F1196").Value = 45824.62041666666
$ws.Range("C1199").Value = 58
$ws.Range("D1199").Value = 45825.28849841943
$ws.Range("E1199").Value = 58
$ws.Range("F1199").Value = 45824.64331018519
$ws.Range("C1223").Value = 15
$ws.Range("D1223").Value = 45825.28849841109
$ws.Range("E1223").Value = 15
$ws.Range("F1223").Value = 45824.37760416666
$ws.Range("C1250").Value = 38
$ws.Range("D1250").Value = 45825.28852326456
$ws.Range("E1250").Value = 38
$ws.Range("F1250").Value = 45824.68135416666
$ws.Range("C1253").Value = 829
$ws.Range("D1253").Value = 45825.28852326369
$ws.Range("E1253").Value = 829
$ws.Range("F1253").Value = 45824.64796296296
$ws.Range("C1330").Value = 11
$ws.Range("D1330").Value = 45825.28849841357
$ws.Range("E1330").Value = 11
$ws.Range("F1330").Value = 45824.43540509259
$ws.Range("C1332").Value = -6
$ws.Range("D1332").Value = 45825.2885232677
$ws.Range("E1332").Value = -6
$ws.Range("F1332").Value = 45824.74089120371
$ws.Range("C1342").Value = 783
$ws.Range("D1342").Value = 45825.28849841825
$ws.Range("E1342").Value = 783
$ws.Range("F1342").Value = 45824.64238425926
$ws.Range("C1390").Value = 47
$ws.Range("D1390").Value = 45825.2884984117
$ws.Range("E1390").Value = 47
$ws.Range("F1390").Value = 45824.3782175926
$ws.Range("C1392").Value = 9
$ws.Range("D1392").Value = 45825.28852326855
$ws.Range("E1392").Value = 9
$ws.Range("F1392").Value = 45824.74362268519
$ws.Range("C1412").Value = 134
$ws.Range("D1412").Value = 45825.28852326781
$ws.Range("E1412").Value = 134
$ws.Range("F1412").Value = 45824.74089120371
$ws.Range("C1446").Value = 54
$ws.Range("D1446").Value = 45825.28849841491
$ws.Range("E1446").Value = 54
$ws.Range("F1446").Value = 45824.45626157407
$ws.Range("C1448").Value = 275
$ws.Range("D1448").Value = 45825.28849841757
$ws.Range("E1448").Value = 275
$ws.Range("F1448").Value = 45824.63096064814
$ws.Range("C1501").Value = 1
$ws.Range("D1501").Value = 45825.2884984112
$ws.Range("E1501").Value = 1
$ws.Range("F1501").Value = 45824.37760416666
$ws.Range("C1507").Value = 85
$ws.Range("D1507").Value = 45825.28852326384
$ws.Range("E1507").Value = 85
$ws.Range("F1507").Value = 45824.64796296296
$ws.Range("C1595").Value = 3
$ws.Range("D1595").Value = 45825.28849841367
$ws.Range("E1595").Value = 3
$ws.Range("F1595").Value = 45824.43540509259
$ws.Range("C1597").Value = 5396
$ws.Range("D1597").Value = 45825.28852326867
$ws.Range("E1597").Value = 5396
$ws.Range("F1597").Value = 45824.74362268519
$ws.Range("C1635").Value = 2671
$ws.Range("D1635").Value = 45825.28849841896
$ws.Range("E1635").Value = 2671
$ws.Range("F1635").Value = 45824.6428125
$ws.Range("C1651").Value = 12
$ws.Range("D1651").Value = 45825.28852326957
$ws.Range("E1651").Value = 12
$ws.Range("F1651").Value = 45824.76962962963
$ws.Range("C1656").Value = 2
$ws.Range("D1656").Value = 45825.28852326066
$ws.Range("E1656").Value = 2
$ws.Range("F1656").Value = 45824.64541666667
$ws.Range("C1844").Value = 256
$ws.Range("D1844").Value = 45825.28849841183
$ws.Range("E1844").Value = 256
$ws.Range("F1844").Value = 45824.3782175926
$ws.Range("C2023").Value = 11
$ws.Range("D2023").Value = 45825.28852326728
$ws.Range("E2023").Value = 11
$ws.Range("F2023").Value = 45824.74057870371
$ws.Range("C2033").Value = 5
$ws.Range("D2033").Value = 45825.28849841766
$ws.Range("E2033").Value = 5
$ws.Range("F2033").Value = 45824.63096064814
$ws.Range("C2045").Value = 202
$ws.Range("D2045").Value = 45825.28849841699
$ws.Range("E2045").Value = 202
$ws.Range("F2045").Value = 45824.62626157407
$ws.Range("C2056").Value = 5
$ws.Range("D2056").Value = 45825.28852326663
$ws.Range("E2056").Value = 5
$ws.Range("F2056").Value = 45824.74033564814
$ws.Range("C2180").Value = 48
$ws.Range("D2180").Value = 45825.28849841524
$ws.Range("E2180").Value = 48
$ws.Range("F2180").Value = 45824.4725
$ws.Range("C2192").Value = 9
$ws.Range("D2192").Value = 45825.28849841432
$ws.Range("E2192").Value = 9
$ws.Range("F2192").Value = 45824.44917824074
$ws.Range("C2240").Value = 268
$ws.Range("D2240").Value = 45825.28852326396
$ws.Range("E2240").Value = 268
$ws.Range("F2240").Value = 45824.64796296296
$ws.Range("C2241").Value = 276
$ws.Range("D2241").Value = 45825.28849841196
$ws.Range("E2241").Value = 276
$ws.Range("F2241").Value = 45824.3782175926
$ws.Range("C2327").Value = 24
$ws.Range("D2327").Value = 45825.2885232648
$ws.Range("E2327").Value = 24
$ws.Range("F2327").Value = 45824.68784722222
$ws.Range("C2328").Value = 5
$ws.Range("D2328").Value = 45825.28852326913
$ws.Range("E2328").Value = 5
$ws.Range("F2328").Value = 45824.75456018518
$ws.Range("C2416").Value = 135
$ws.Range("D2416").Value = 45825.28849842033
$ws.Range("E2416").Value = 135
$ws.Range("F2416").Value = 45824.64493055556
$ws.Range("C2481").Value = 408
$ws.Range("D2481").Value = 45825.28849841421
$ws.Range("E2481").Value = 408
$ws.Range("F2481").Value = 45824.44886574074
$ws.Range("C2483").Value = 463
$ws.Range("D2483").Value = 45825.28849841573
$ws.Range("E2483").Value = 463
$ws.Range("F2483").Value = 45824.50672453704

# Row 2570 previously had no F (data_movimento) value; add it with the
# same date/time number format used by the rest of column F/D.
$ws.Range("F2570").Value = 45824.64381944444
$ws.Range("F2570").NumberFormat = $ws.Range("D2570").NumberFormat

# Three brand-new product rows appended at the end of the sheet, with the
# same column layout / status text as the existing consistent rows (no
# movement recorded yet, hence no F value).
$ws.Range("A2573").Value = 43660986
$ws.Range("B2573").Value = 1
$ws.Range("C2573").Value = 0
$ws.Range("D2573").Value = 45825.28852326623
$ws.Range("D2573").NumberFormat = $ws.Range("D2572").NumberFormat
$ws.Range("E2573").Value = 0
$ws.Range("G2573").Value = 0
$ws.Range("H2573").Value = "Consistente"

$ws.Range("A2574").Value = 43660992
$ws.Range("B2574").Value = 1
$ws.Range("C2574").Value = 0
$ws.Range("D2574").Value = 45825.28852326625
$ws.Range("D2574").NumberFormat = $ws.Range("D2572").NumberFormat
$ws.Range("E2574").Value = 0
$ws.Range("G2574").Value = 0
$ws.Range("H2574").Value = "Consistente"

$ws.Range("A2575").Value = 43664584
$ws.Range("B2575").Value = 1
$ws.Range("C2575").Value = 0
$ws.Range("D2575").Value = 45825.2885232687
$ws.Range("D2575").NumberFormat = $ws.Range("D2572").NumberFormat
$ws.Range("E2575").Value = 0
$ws.Range("G2575").Value = 0
$ws.Range("H2575").Value = "Consistente"
